$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata" (sheet1): StructureDefinition summary table ---
$ws1 = $wb.Worksheets.Item("Metadata")

# Remove the old standalone "Description" row (previously row 12); its content
# (Description / Numeric version of the organinzational hierarchy represented)
# is being relocated into what is currently the second "Contact" row.
$ws1.Rows.Item(12).Delete()

# Bump the published version number.
$ws1.Range("B3").Value = "6.0.0"

# Update the publication date.
$ws1.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher is now populated.
$ws1.Range("B9").Value = "Alvearie Team"

# First "Contact" row becomes "Jurisdiction".
$ws1.Range("A10").Value = "Jurisdiction"
$ws1.Range("B10").Value = "United States of America"

# Second "Contact" row becomes "Description".
$ws1.Range("A11").Value = "Description"
$ws1.Range("B11").Value = "Numeric version of the organinzational hierarchy represented"

# --- Sheet "Elements" (sheet2): element definitions table ---
$ws2 = $wb.Worksheets.Item("Elements")

# Root element row ("Extension") gets the resource-specific Short / Definition text
# instead of the generic placeholders.
$ws2.Range("K2").Value = "Hierarchy VersionId"
$ws2.Range("L2").Value = "Numeric version of the organinzational hierarchy represented"
